# Fix formatting when scraping floating point numbers:
# - Column H (Importe) values were stored as Spanish/AR-formatted
#   text ("." thousands, "," decimal); convert to plain decimal
#   text with "." as the decimal separator, keeping them as text.
# - A handful of Razon social (column E) entries had stray commas
# that should read as periods.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H holds these as text (General format) in the source file;
# force Text format so plain-looking numeric strings (e.g. "11730.00")
# are not silently reinterpreted as numbers and lose their formatting.
$ws.Range("H2:H279").NumberFormat = "@"

$importeValues = @(
    "11730.00",
    "40960.00",
    "4050.00",
    "440000.00",
    "195000.00",
    "561000.00",
    "665500.00",
    "2250.00",
    "545773.47",
    "582852.20",
    "159281.54",
    "181500.00",
    "4000.00",
    "10000.00",
    "6500.06",
    "360.00",
    "3411.00",
    "770.00",
    "88000.00",
    "490203.00",
    "232584.81",
    "25600.00",
    "840.00",
    "346677.74",
    "9765.00",
    "3120.00",
    "526.10",
    "59125.53",
    "26503.00",
    "26103.22",
    "20600.00",
    "9450.00",
    "9200.00",
    "21000.00",
    "756.62",
    "116.00",
    "7045.00",
    "422.25",
    "234269.66",
    "365.00",
    "134201.00",
    "226500.00",
    "510.00",
    "5818.49",
    "860.96",
    "1898.00",
    "1246.00",
    "22148.72",
    "3080.00",
    "5929.00",
    "2200.00",
    "7098.20",
    "85.00",
    "21546.27",
    "8970.00",
    "6296.10",
    "8113.10",
    "179676.10",
    "2627.39",
    "34548.34",
    "436488.80",
    "1660.00",
    "4500.00",
    "240.00",
    "65.00",
    "60838.12",
    "358577.00",
    "4740.00",
    "4900.00",
    "11663.85",
    "16980.00",
    "2200.00",
    "6900.00",
    "9450.00",
    "1527.00",
    "3793.62",
    "4701.86",
    "2600.00",
    "151.05",
    "7611.79",
    "5205.90",
    "26950.00",
    "7250.00",
    "1400.00",
    "9352.50",
    "490.00",
    "350.00",
    "120.00",
    "92070.00",
    "2974.00",
    "20625.00",
    "8690.00",
    "14920.00",
    "6555.00",
    "482.00",
    "39000.00",
    "45900.00",
    "4050.00",
    "2352.00",
    "4070.00",
    "17.32",
    "11.20",
    "21000.00",
    "330.55",
    "28.90",
    "100.00",
    "28596.67",
    "936.00",
    "1335.05",
    "231.45",
    "1254.00",
    "53.80",
    "2400.00",
    "11778.80",
    "32090.00",
    "21864.00",
    "285.00",
    "1879.00",
    "2298.50",
    "10050.00",
    "1256.54",
    "2430.00",
    "129.60",
    "9579.00",
    "7000.00",
    "15123.50",
    "32572.21",
    "670.10",
    "2541.50",
    "100.35",
    "4430.00",
    "2012.90",
    "8760.00",
    "4700.00",
    "4626.00",
    "7500.00",
    "500.00",
    "7194.00",
    "160000.00",
    "2505.00",
    "54700.00",
    "7618.00",
    "3500.00",
    "72000.00",
    "75200.00",
    "14000.00",
    "12000.00",
    "44000.00",
    "3500.00",
    "36000.00",
    "60000.00",
    "3689.41",
    "11675.00",
    "16934.21",
    "6748.88",
    "4889.00",
    "200.00",
    "7784.23",
    "25000.00",
    "12500.00",
    "2000.00",
    "14000.00",
    "8000.00",
    "52111.68",
    "4500.00",
    "5000.00",
    "8508.50",
    "10500.00",
    "5000.00",
    "6000.00",
    "6000.00",
    "5000.00",
    "5000.00",
    "6000.00",
    "7000.00",
    "23000.00",
    "6000.00",
    "20000.00",
    "15000.00",
    "11500.00",
    "6000.00",
    "3500.00",
    "45000.00",
    "10000.00",
    "100672.00",
    "24500.00",
    "9800.00",
    "9800.00",
    "67000.00",
    "23050.00",
    "5293.19",
    "574.07",
    "3814.00",
    "4275.00",
    "3337.00",
    "37617.87",
    "13738.00",
    "8262.00",
    "36.00",
    "22445.50",
    "5715.00",
    "1210.00",
    "7199.15",
    "655.74",
    "117000.00",
    "4781.20",
    "28166.04",
    "1710.00",
    "21319.98",
    "28590.00",
    "4269.12",
    "16004.90",
    "360.00",
    "7382.00",
    "110.00",
    "1092.60",
    "5674.00",
    "467.50",
    "1304.00",
    "16272.00",
    "7306.24",
    "21534.74",
    "8500.00",
    "30000.00",
    "30000.00",
    "30000.00",
    "30000.00",
    "30000.00",
    "30000.00",
    "60000.00",
    "60000.00",
    "30000.00",
    "9000.00",
    "1148.25",
    "7000.00",
    "14400.00",
    "2390.00",
    "13950.00",
    "3850698.08",
    "135500.00",
    "144500.00",
    "135500.00",
    "145000.00",
    "135500.00",
    "135500.00",
    "239000.00",
    "324000.00",
    "342500.00",
    "135500.00",
    "135500.00",
    "135500.00",
    "135500.00",
    "135500.00",
    "135500.00",
    "342500.00",
    "135500.00",
    "244000.00",
    "135500.00",
    "135500.00",
    "4750.00",
    "135500.00",
    "54667.53",
    "5345897.56",
    "80000.00",
    "9000.00",
    "35040.55",
    "6700.00",
    "16000.00",
    "117000.00",
    "36196.66",
    "60000.00",
    "70000.00",
    "13499.99",
    "8550.00",
    "27072.50",
    "3000.00",
    "480.00",
    "26168.00"
)

for ($i = 0; $i -lt $importeValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $importeValues[$i]
}

# Razon social punctuation fixes (comma -> period; drop stray dots in "S.H.")
$ws.Cells.Item(41, 5).Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Cells.Item(50, 5).Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
$ws.Cells.Item(92, 5).Value = "FERNANDEZ. MARIO HUGO"
$ws.Cells.Item(94, 5).Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Cells.Item(194, 5).Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
$ws.Cells.Item(211, 5).Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"

